# Edit derived from the OOXML diff for xl/worksheets/sheet1.xml
# Summary of the change:
#   - Sheet grows from A1:O25 to A1:Q25 (two new columns: P, Q)
#   - Row 1 (header) gains P1=14, Q1=15, styled like the existing header cells
#   - Rows 2-25 get fully refreshed numeric data (B:I recomputed, O zeroed out,
#     P added as 0, and the former "last column" value now lives in the new Q column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add P1 / Q1, copying O1's format (bold/centered/bordered style) ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)   # xlPasteFormats: style only, no value
$ws.Application.CutCopyMode = $false
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# --- Data rows 2-25, columns A:Q -- bulk write via a real 2-D array ---
$arr = New-Object 'double[,]' 24,17
$arr[0,0] = 0
$arr[0,1] = 3.46677384630982
$arr[0,2] = 1.019169272890139
$arr[0,3] = 0.04824190763466873
$arr[0,4] = 1.35169962688731
$arr[0,5] = 0.5084535737087279
$arr[0,6] = 0.0007896385517037115
$arr[0,7] = 0.01092094002687105
$arr[0,8] = 0.003536527785124033
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0
$arr[0,15] = 0
$arr[0,16] = 1.435151477765999
$arr[1,0] = 1
$arr[1,1] = 3.02011174815118
$arr[1,2] = 0.8990605546520669
$arr[1,3] = 0.0433835340809452
$arr[1,4] = 1.176740016824183
$arr[1,5] = 0.468187038062311
$arr[1,6] = 0.000793366441376707
$arr[1,7] = 0.00766975902495437
$arr[1,8] = 0.001965288275290966
$arr[1,9] = 0
$arr[1,10] = 0
$arr[1,11] = 0
$arr[1,12] = 0
$arr[1,13] = 0
$arr[1,14] = 0
$arr[1,15] = 0
$arr[1,16] = 1.348594257482659
$arr[2,0] = 2
$arr[2,1] = 2.745669458445775
$arr[2,2] = 0.8258832939075091
$arr[2,3] = 0.04041041276209967
$arr[2,4] = 1.069666178277558
$arr[2,5] = 0.4441198563540496
$arr[2,6] = 0.0007957272508090085
$arr[2,7] = 0.005904022270549425
$arr[2,8] = 0.001270814676121557
$arr[2,9] = 0
$arr[2,10] = 0
$arr[2,11] = 0
$arr[2,12] = 0
$arr[2,13] = 0
$arr[2,14] = 0
$arr[2,15] = 0
$arr[2,16] = 1.297464469391542
$arr[3,0] = 3
$arr[3,1] = 2.633730900298815
$arr[3,2] = 0.7977722387358028
$arr[3,3] = 0.03928077046956702
$arr[3,4] = 1.02608314129057
$arr[3,5] = 0.4336700664951891
$arr[3,6] = 0.0007967124296093884
$arr[3,7] = 0.00523637476029204
$arr[3,8] = 0.001110499418202515
$arr[3,9] = 0
$arr[3,10] = 0
$arr[3,11] = 0
$arr[3,12] = 0
$arr[3,13] = 0
$arr[3,14] = 0
$arr[3,15] = 0
$arr[3,16] = 1.274436995028537
$arr[4,0] = 4
$arr[4,1] = 2.615103766652396
$arr[4,2] = 0.7950489896263662
$arr[4,3] = 0.03919122151063448
$arr[4,4] = 1.018830727099584
$arr[4,5] = 0.4309748955527013
$arr[4,6] = 0.0007968828780179578
$arr[4,7] = 0.005126972417599962
$arr[4,8] = 0.001168165231651308
$arr[4,9] = 0
$arr[4,10] = 0
$arr[4,11] = 0
$arr[4,12] = 0
$arr[4,13] = 0
$arr[4,14] = 0
$arr[4,15] = 0
$arr[4,16] = 1.267397047251904
$arr[5,0] = 5
$arr[5,1] = 2.74406360535454
$arr[5,2] = 0.8308092758339853
$arr[5,3] = 0.04066341026042863
$arr[5,4] = 1.069026406947373
$arr[5,5] = 0.4413300423773592
$arr[5,6] = 0.0007957557770997617
$arr[5,7] = 0.005890106154369357
$arr[5,8] = 0.001472853056091239
$arr[5,9] = 0
$arr[5,10] = 0
$arr[5,11] = 0
$arr[5,12] = 0
$arr[5,13] = 0
$arr[5,14] = 0
$arr[5,15] = 0
$arr[5,16] = 1.28828808501018
$arr[6,0] = 6
$arr[6,1] = 3.312653545137607
$arr[6,2] = 0.9847298999007421
$arr[6,3] = 0.04692401796475565
$arr[6,4] = 1.291216389130142
$arr[6,5] = 0.4909136236972671
$arr[6,6] = 0.0007909279253908628
$arr[6,7] = 0.009741767389663553
$arr[6,8] = 0.003174773938467546
$arr[6,9] = 0
$arr[6,10] = 0
$arr[6,11] = 0
$arr[6,12] = 0
$arr[6,13] = 0
$arr[6,14] = 0
$arr[6,15] = 0
$arr[6,16] = 1.393125671300481
$arr[7,0] = 7
$arr[7,1] = 4.427610317638539
$arr[7,2] = 1.282110248089054
$arr[7,3] = 0.0588090931990024
$arr[7,4] = 1.730915313986756
$arr[7,5] = 0.5988199823157387
$arr[7,6] = 0.0007819898167375172
$arr[7,7] = 0.01935802986678037
$arr[7,8] = 0.008607690668569923
$arr[7,9] = 0
$arr[7,10] = 0
$arr[7,11] = 0
$arr[7,12] = 0
$arr[7,13] = 0
$arr[7,14] = 0
$arr[7,15] = 0
$arr[7,16] = 1.632718893241275
$arr[8,0] = 8
$arr[8,1] = 5.24749449069725
$arr[8,2] = 1.50644471872306
$arr[8,3] = 0.06925275944912102
$arr[8,4] = 1.955498156088424
$arr[8,5] = 0.6683724546327312
$arr[8,6] = 0.0007759177289413412
$arr[8,7] = 0.02725836382229918
$arr[8,8] = 0.0143664161419883
$arr[8,9] = 0
$arr[8,10] = 0
$arr[8,11] = 0
$arr[8,12] = 0
$arr[8,13] = 0
$arr[8,14] = 0
$arr[8,15] = 0
$arr[8,16] = 1.775547518866148
$arr[9,0] = 9
$arr[9,1] = 5.612142301903873
$arr[9,2] = 1.610434802784823
$arr[9,3] = 0.08861588644301577
$arr[9,4] = 1.267337406619561
$arr[9,5] = 0.5871507068896307
$arr[9,6] = 0.0007746196692268896
$arr[9,7] = 0.04270868006982909
$arr[9,8] = 0.01643561128637661
$arr[9,9] = 0
$arr[9,10] = 0
$arr[9,11] = 0
$arr[9,12] = 0
$arr[9,13] = 0
$arr[9,14] = 0
$arr[9,15] = 0
$arr[9,16] = 1.474547141099578
$arr[10,0] = 10
$arr[10,1] = 5.74725555498668
$arr[10,2] = 1.642389293035194
$arr[10,3] = 0.1037805964393641
$arr[10,4] = 0.7713189373199754
$arr[10,5] = 0.5121260620343264
$arr[10,6] = 0.0007745905169331831
$arr[10,7] = 0.07854308646568597
$arr[10,8] = 0.01664005842228544
$arr[10,9] = 0
$arr[10,10] = 0
$arr[10,11] = 0
$arr[10,12] = 0
$arr[10,13] = 0
$arr[10,14] = 0
$arr[10,15] = 0
$arr[10,16] = 1.224683832186997
$arr[11,0] = 11
$arr[11,1] = 5.71206086016997
$arr[11,2] = 1.628934670659646
$arr[11,3] = 0.1168170525753141
$arr[11,4] = 0.3948377359290731
$arr[11,5] = 0.4334589186537912
$arr[11,6] = 0.0007755528712741501
$arr[11,7] = 0.1314041879758321
$arr[11,8] = 0.01567717074598907
$arr[11,9] = 0
$arr[11,10] = 0
$arr[11,11] = 0
$arr[11,12] = 0
$arr[11,13] = 0
$arr[11,14] = 0
$arr[11,15] = 0
$arr[11,16] = 0.9858811040282944
$arr[12,0] = 12
$arr[12,1] = 5.610402796626545
$arr[12,2] = 1.600782500645209
$arr[12,3] = 0.1251548963565625
$arr[12,4] = 0.2040035724406977
$arr[12,5] = 0.3772055239524477
$arr[12,6] = 0.0007766576504858142
$arr[12,7] = 0.1788295655706946
$arr[12,8] = 0.01459786682350295
$arr[12,9] = 0
$arr[12,10] = 0
$arr[12,11] = 0
$arr[12,12] = 0
$arr[12,13] = 0
$arr[12,14] = 0
$arr[12,15] = 0
$arr[12,16] = 0.8260307867355579
$arr[13,0] = 13
$arr[13,1] = 5.549222564864976
$arr[13,2] = 1.586479300100507
$arr[13,3] = 0.1266680095887693
$arr[13,4] = 0.1664673984577796
$arr[13,5] = 0.361549071281118
$arr[13,6] = 0.0007771664842525621
$arr[13,7] = 0.1906857109769646
$arr[13,8] = 0.01418120435697556
$arr[13,9] = 0
$arr[13,10] = 0
$arr[13,11] = 0
$arr[13,12] = 0
$arr[13,13] = 0
$arr[13,14] = 0
$arr[13,15] = 0
$arr[13,16] = 0.7848470606577962
$arr[14,0] = 14
$arr[14,1] = 5.202401853822494
$arr[14,2] = 1.495615476278033
$arr[14,3] = 0.1191748546277722
$arr[14,4] = 0.1619581450697751
$arr[14,5] = 0.3473887037929586
$arr[14,6] = 0.0007795445118418047
$arr[14,7] = 0.175491235802042
$arr[14,8] = 0.01192358401845794
$arr[14,9] = 0
$arr[14,10] = 0
$arr[14,11] = 0
$arr[14,12] = 0
$arr[14,13] = 0
$arr[14,14] = 0
$arr[14,15] = 0
$arr[14,16] = 0.7706965671103063
$arr[15,0] = 15
$arr[15,1] = 4.991035451832488
$arr[15,2] = 1.441267988981963
$arr[15,3] = 0.1089952125854552
$arr[15,4] = 0.2459753982704243
$arr[15,5] = 0.365796307107324
$arr[15,6] = 0.000780802229109395
$arr[15,7] = 0.1367689001963726
$arr[15,8] = 0.01081278822335374
$arr[15,9] = 0
$arr[15,10] = 0
$arr[15,11] = 0
$arr[15,12] = 0
$arr[15,13] = 0
$arr[15,14] = 0
$arr[15,15] = 0
$arr[15,16] = 0.8422055574621652
$arr[16,0] = 16
$arr[16,1] = 4.871728470037112
$arr[16,2] = 1.407347280808153
$arr[16,3] = 0.09563455489917061
$arr[16,4] = 0.4790218023010837
$arr[16,5] = 0.4181691898324047
$arr[16,6] = 0.0007811315644281613
$arr[16,7] = 0.08461431799794639
$arr[16,8] = 0.01026415122278657
$arr[16,9] = 0
$arr[16,10] = 0
$arr[16,11] = 0
$arr[16,12] = 0
$arr[16,13] = 0
$arr[16,14] = 0
$arr[16,15] = 0
$arr[16,16] = 1.011170216096644
$arr[17,0] = 17
$arr[17,1] = 4.83500880465499
$arr[17,2] = 1.40339474351282
$arr[17,3] = 0.08223365033970964
$arr[17,4] = 0.9131096264680565
$arr[17,5] = 0.4940337535870114
$arr[17,6] = 0.0007805737522911431
$arr[17,7] = 0.0420758044946794
$arr[17,8] = 0.01072812603060136
$arr[17,9] = 0
$arr[17,10] = 0
$arr[17,11] = 0
$arr[17,12] = 0
$arr[17,13] = 0
$arr[17,14] = 0
$arr[17,15] = 0
$arr[17,16] = 1.251779364725252
$arr[18,0] = 18
$arr[18,1] = 5.031627643970751
$arr[18,2] = 1.464437765624723
$arr[18,3] = 0.06744620694644254
$arr[18,4] = 1.891890550245648
$arr[18,5] = 0.6409981226065185
$arr[18,6] = 0.0007775358304290444
$arr[18,7] = 0.02500230937025849
$arr[18,8] = 0.01333659226125761
$arr[18,9] = 0
$arr[18,10] = 0
$arr[18,11] = 0
$arr[18,12] = 0
$arr[18,13] = 0
$arr[18,14] = 0
$arr[18,15] = 0
$arr[18,16] = 1.707873429777266
$arr[19,0] = 19
$arr[19,1] = 5.663603663768583
$arr[19,2] = 1.636139845533137
$arr[19,3] = 0.07300782904698622
$arr[19,4] = 2.224275313971731
$arr[19,5] = 0.7181090159483858
$arr[19,6] = 0.0007726865873697398
$arr[19,7] = 0.03281830044685519
$arr[19,8] = 0.01854060229383236
$arr[19,9] = 0
$arr[19,10] = 0
$arr[19,11] = 0
$arr[19,12] = 0
$arr[19,13] = 0
$arr[19,14] = 0
$arr[19,15] = 0
$arr[19,16] = 1.895197017569245
$arr[20,0] = 20
$arr[20,1] = 6.077296100170315
$arr[20,2] = 1.741715857843474
$arr[20,3] = 0.07705961331203781
$arr[20,4] = 2.390794747381719
$arr[20,5] = 0.7665380570823004
$arr[20,6] = 0.000769642109448282
$arr[20,7] = 0.03795447932891705
$arr[20,8] = 0.02207615676014552
$arr[20,9] = 0
$arr[20,10] = 0
$arr[20,11] = 0
$arr[20,12] = 0
$arr[20,13] = 0
$arr[20,14] = 0
$arr[20,15] = 0
$arr[20,16] = 2.010982655735347
$arr[21,0] = 21
$arr[21,1] = 5.856533976777541
$arr[21,2] = 1.678693843975964
$arr[21,3] = 0.0745603636229859
$arr[21,4] = 2.301860511110718
$arr[21,5] = 0.7437029678328884
$arr[21,6] = 0.0007712498272074986
$arr[21,7] = 0.03518443906473134
$arr[21,8] = 0.01990882957891671
$arr[21,9] = 0
$arr[21,10] = 0
$arr[21,11] = 0
$arr[21,12] = 0
$arr[21,13] = 0
$arr[21,14] = 0
$arr[21,15] = 0
$arr[21,16] = 1.959238721746175
$arr[22,0] = 22
$arr[22,1] = 5.022345505370311
$arr[22,2] = 1.452387582760707
$arr[22,3] = 0.0656626206233426
$arr[22,4] = 1.967538615152051
$arr[22,5] = 0.6548386594035378
$arr[22,6] = 0.0007774558731151467
$arr[22,7] = 0.02548343545092102
$arr[22,8] = 0.01297716825878936
$arr[22,9] = 0
$arr[22,10] = 0
$arr[22,11] = 0
$arr[22,12] = 0
$arr[22,13] = 0
$arr[22,14] = 0
$arr[22,15] = 0
$arr[22,16] = 1.754635322995938
$arr[23,0] = 23
$arr[23,1] = 4.125615602166818
$arr[23,2] = 1.210872665251486
$arr[23,3] = 0.05608168448596018
$arr[23,4] = 1.611287029704712
$arr[23,5] = 0.5640920550256467
$arr[23,6] = 0.0007843738772956623
$arr[23,7] = 0.01649207729639268
$arr[23,8] = 0.007186637996994705
$arr[23,9] = 0
$arr[23,10] = 0
$arr[23,11] = 0
$arr[23,12] = 0
$arr[23,13] = 0
$arr[23,14] = 0
$arr[23,15] = 0
$arr[23,16] = 1.549647283128024

$ws.Range("A2:Q25").Value2 = $arr

